# Auto-generated cell updates applying the Odin_Profits.xlsx diff
# (scheduled-runner style bulk value overwrite; no formulas involved)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 4499.091
$ws.Range("J33").Value = 1025
$ws.Range("L33").Value = 1025
$ws.Range("N33").Value = -1483
# Row 39
$ws.Range("H39").Value = 792.63635
$ws.Range("I39").Value = 80
$ws.Range("J39").Value = 3999.5
$ws.Range("K39").Value = 240
$ws.Range("L39").Value = 11998.5
$ws.Range("M39").Value = 56
$ws.Range("N39").Value = -12590.5
# Row 40
$ws.Range("H40").Value = 2694.75
$ws.Range("J40").Value = 2399
$ws.Range("L40").Value = 2399
$ws.Range("N40").Value = -2749
# Row 41
$ws.Range("H41").Value = 443.33334
$ws.Range("I41").Value = 443.33334
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 443.33334
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -3.333340000000021
# Row 42
$ws.Range("H42").Value = 200
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
# Row 82
$ws.Range("H82").Value = 2082.5715
$ws.Range("I82").Value = 755.1667
$ws.Range("J82").Value = 10047
$ws.Range("K82").Value = 2265.5001
$ws.Range("L82").Value = 30141
$ws.Range("M82").Value = -1859.5001
$ws.Range("N82").Value = -30953
# Row 85
$ws.Range("H85").Value = 2082.5715
$ws.Range("I85").Value = 755.1667
$ws.Range("J85").Value = 10047
$ws.Range("K85").Value = 2265.5001
$ws.Range("L85").Value = 30141
$ws.Range("M85").Value = -861.5001000000002
$ws.Range("N85").Value = -32949
# Row 96
$ws.Range("H96").Value = 730.2
$ws.Range("I96").Value = 579.625
$ws.Range("J96").Value = 1332.5
$ws.Range("K96").Value = 1738.875
$ws.Range("L96").Value = 3997.5
$ws.Range("M96").Value = -365.875
$ws.Range("N96").Value = -6743.5
# Row 98
$ws.Range("H98").Value = 1700.7354
$ws.Range("I98").Value = 954.931
$ws.Range("K98").Value = 954.931
$ws.Range("M98").Value = 543.069
# Row 101
$ws.Range("H101").Value = 1356.2667
$ws.Range("I101").Value = 1537.1
$ws.Range("K101").Value = 4611.299999999999
$ws.Range("M101").Value = -2989.299999999999
# Row 122
$ws.Range("H122").Value = 1700.7354
$ws.Range("I122").Value = 954.931
$ws.Range("K122").Value = 2864.793
$ws.Range("M122").Value = -414.7930000000001
# Row 131
$ws.Range("H131").Value = 1427.125
$ws.Range("I131").Value = 1059.5714
$ws.Range("K131").Value = 3178.7142
$ws.Range("M131").Value = 1861.2858
# Row 138
$ws.Range("H138").Value = 4116.4287
$ws.Range("J138").Value = 4155.3726
$ws.Range("L138").Value = 12466.1178
$ws.Range("N138").Value = -22746.1178

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 5469.905
$ws.Range("I61").Value = 7614.0713
$ws.Range("K61").Value = 7614.0713
$ws.Range("M61").Value = -7402.0713
# Row 63
$ws.Range("H63").Value = 2209.389
$ws.Range("J63").Value = 2964.6
$ws.Range("L63").Value = 2964.6
$ws.Range("N63").Value = -4336.6
# Row 66
$ws.Range("H66").Value = 2209.389
$ws.Range("J66").Value = 2964.6
$ws.Range("L66").Value = 14823
$ws.Range("N66").Value = -21687
# Row 97
$ws.Range("H97").Value = 512.6842
$ws.Range("J97").Value = 617.8333
$ws.Range("L97").Value = 617.8333
$ws.Range("N97").Value = -1609.8333
# Row 136
$ws.Range("H136").Value = 5469.905
$ws.Range("I136").Value = 7614.0713
$ws.Range("K136").Value = 22842.2139
$ws.Range("M136").Value = -20292.2139

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 42741.69
$ws.Range("I94").Value = 1182.5385
$ws.Range("J94").Value = 84300.84
$ws.Range("K94").Value = 1182.5385
$ws.Range("L94").Value = 84300.84
$ws.Range("M94").Value = -731.5385000000001
$ws.Range("N94").Value = -85202.84

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 26320006
$ws.Range("J16").Value = 6251.25
$ws.Range("L16").Value = 6251.25
$ws.Range("N16").Value = -6825.25
# Row 31
$ws.Range("H31").Value = 2632.3062
$ws.Range("I31").Value = 913.3946999999999
$ws.Range("J31").Value = 8570.362999999999
$ws.Range("K31").Value = 913.3946999999999
$ws.Range("L31").Value = 8570.362999999999
$ws.Range("M31").Value = -618.3946999999999
$ws.Range("N31").Value = -9160.362999999999
# Row 34
$ws.Range("H34").Value = 2632.3062
$ws.Range("I34").Value = 913.3946999999999
$ws.Range("J34").Value = 8570.362999999999
$ws.Range("K34").Value = 913.3946999999999
$ws.Range("L34").Value = 8570.362999999999
$ws.Range("M34").Value = -711.3946999999999
$ws.Range("N34").Value = -8974.362999999999
# Row 62
$ws.Range("H62").Value = 5286.905
$ws.Range("J62").Value = 5559.8335
$ws.Range("L62").Value = 5559.8335
$ws.Range("N62").Value = -6807.8335
# Row 65
$ws.Range("H65").Value = 5286.905
$ws.Range("J65").Value = 5559.8335
$ws.Range("L65").Value = 27799.1675
$ws.Range("N65").Value = -34039.1675
# Row 113
$ws.Range("H113").Value = 26320006
$ws.Range("J113").Value = 6251.25
$ws.Range("L113").Value = 6251.25
$ws.Range("N113").Value = -10591.25
# Row 140
$ws.Range("H140").Value = 99999
$ws.Range("J140").Value = 129998
$ws.Range("L140").Value = 129998
$ws.Range("N140").Value = -140358

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 74
$ws.Range("H74").Value = 21874.5
$ws.Range("J74").Value = 28749.5
$ws.Range("L74").Value = 86248.5
$ws.Range("N74").Value = -88370.5
# Row 77
$ws.Range("H77").Value = 21874.5
$ws.Range("J77").Value = 28749.5
$ws.Range("L77").Value = 258745.5
$ws.Range("N77").Value = -269353.5
# Row 114
$ws.Range("H114").Value = 126084.75
$ws.Range("I114").Value = 884.5
$ws.Range("J114").Value = 251285
$ws.Range("K114").Value = 2653.5
$ws.Range("L114").Value = 753855
$ws.Range("M114").Value = 600.5
$ws.Range("N114").Value = -760363
# Row 118
$ws.Range("H118").Value = 16357.25
$ws.Range("I118").Value = 24029
$ws.Range("K118").Value = 72087
$ws.Range("M118").Value = -70844

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 9566.625
$ws.Range("I80").Value = 3965.3333
$ws.Range("K80").Value = 3965.3333
$ws.Range("M80").Value = -2967.3333
# Row 83
$ws.Range("H83").Value = 9566.625
$ws.Range("I83").Value = 3965.3333
$ws.Range("K83").Value = 19826.6665
$ws.Range("M83").Value = -14834.6665
# Row 132
$ws.Range("H132").Value = 30318826
$ws.Range("I132").Value = 33350422
$ws.Range("J132").Value = 2865.6667
$ws.Range("K132").Value = 100051266
$ws.Range("L132").Value = 8597.000100000001
$ws.Range("M132").Value = -100048736
$ws.Range("N132").Value = -13657.0001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 29414968
$ws.Range("J46").Value = 50004908
$ws.Range("L46").Value = 50004908
$ws.Range("N46").Value = -50005284
# Row 82
$ws.Range("H82").Value = 2279.3333
$ws.Range("I82").Value = 1619.6666
$ws.Range("J82").Value = 3103.9167
$ws.Range("K82").Value = 1619.6666
$ws.Range("L82").Value = 3103.9167
$ws.Range("M82").Value = -1258.6666
$ws.Range("N82").Value = -3825.9167
# Row 85
$ws.Range("H85").Value = 2279.3333
$ws.Range("I85").Value = 1619.6666
$ws.Range("J85").Value = 3103.9167
$ws.Range("K85").Value = 1619.6666
$ws.Range("L85").Value = 3103.9167
$ws.Range("M85").Value = -371.6666
$ws.Range("N85").Value = -5599.9167
# Row 93
$ws.Range("H93").Value = 805.4
$ws.Range("I93").Value = 756.375
$ws.Range("K93").Value = 756.375
$ws.Range("M93").Value = 491.625
# Row 132
$ws.Range("H132").Value = 4806.355
$ws.Range("I132").Value = 4976.92
$ws.Range("K132").Value = 14930.76
$ws.Range("M132").Value = -12400.76

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 31000.25
$ws.Range("I4").Value = 31000.25
$ws.Range("K4").Value = 31000.25
$ws.Range("M4").Value = -30887.25
# Row 64
$ws.Range("H64").Value = 72025.75
$ws.Range("I64").Value = 40103
$ws.Range("J64").Value = 82666.664
$ws.Range("K64").Value = 40103
$ws.Range("L64").Value = 82666.664
$ws.Range("M64").Value = -39855
$ws.Range("N64").Value = -83162.664
# Row 67
$ws.Range("H67").Value = 72025.75
$ws.Range("I67").Value = 40103
$ws.Range("J67").Value = 82666.664
$ws.Range("K67").Value = 40103
$ws.Range("L67").Value = 82666.664
$ws.Range("M67").Value = -39245
$ws.Range("N67").Value = -84382.664
# Row 100
$ws.Range("H100").Value = 1217.4445
$ws.Range("I100").Value = 1125.6666
$ws.Range("K100").Value = 2251.3332
$ws.Range("M100").Value = -1710.3332
# Row 113
$ws.Range("H113").Value = 5753296
$ws.Range("J113").Value = 5350.385
$ws.Range("L113").Value = 16051.155
$ws.Range("N113").Value = -20391.155
# Row 132
$ws.Range("H132").Value = 5672.067
$ws.Range("I132").Value = 4548.9443
$ws.Range("K132").Value = 13646.8329
$ws.Range("M132").Value = -11116.8329

